# Updates cryptos list (Price / Volume(1h) columns, plus a couple of
# reordered rows) to match the refreshed data pulled on
# Fri Jun 21 17:36:02 UTC 2024.
#
# Note: several "Price" values are plain decimal numbers (e.g. "1.00",
# "0.999", "7.08"). Assigning those as plain strings would let Excel's
# COM layer auto-convert them into numeric cell values, which would
# silently strip the original text formatting (trailing zeros, etc.).
# To keep them as literal text - matching the original inlineStr cells -
# we prefix those values with a leading apostrophe, which is Excel's
# native "treat this as text" quote-prefix convention.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '63.595.92'
$ws.Cells.Item(2, 5).Value = '  -1.79%  '
$ws.Cells.Item(3, 4).Value = '3.477.34'
$ws.Cells.Item(3, 5).Value = '  -1.19%  '
$ws.Cells.Item(4, 5).Value = '  +0.02%  '
$ws.Cells.Item(5, 4).Value = '''580.92'
$ws.Cells.Item(5, 5).Value = '  -2.55%  '
$ws.Cells.Item(6, 4).Value = '''129.53'
$ws.Cells.Item(6, 5).Value = '  -3.37%  '
$ws.Cells.Item(7, 4).Value = '3.480.68'
$ws.Cells.Item(9, 5).Value = '  -2.00%  '
$ws.Cells.Item(11, 4).Value = '''7.08'
$ws.Cells.Item(11, 5).Value = '  -0.97%  '
$ws.Cells.Item(12, 5).Value = '  -1.94%  '
$ws.Cells.Item(13, 4).Value = '4.069.25'
$ws.Cells.Item(13, 5).Value = '  -1.08%  '
$ws.Cells.Item(14, 4).Value = '''27.05'
$ws.Cells.Item(14, 5).Value = '  -2.12%  '
$ws.Cells.Item(15, 5).Value = '  +1.10%  '
$ws.Cells.Item(16, 4).Value = '3.489.81'
$ws.Cells.Item(16, 5).Value = '  -0.75%  '
$ws.Cells.Item(17, 4).Value = '''0.0000175'
$ws.Cells.Item(17, 5).Value = '  -3.75%  '
$ws.Cells.Item(18, 4).Value = '63.667.87'
$ws.Cells.Item(18, 5).Value = '  -1.82%  '
$ws.Cells.Item(19, 4).Value = '''9.78'
$ws.Cells.Item(19, 5).Value = '  -3.06%  '
$ws.Cells.Item(20, 4).Value = '''14.01'
$ws.Cells.Item(20, 5).Value = '  -2.29%  '
$ws.Cells.Item(21, 5).Value = '  -1.90%  '
$ws.Cells.Item(22, 4).Value = '''379.22'
$ws.Cells.Item(22, 5).Value = '  -3.07%  '
$ws.Cells.Item(23, 5).Value = '  -1.45%  '
$ws.Cells.Item(24, 4).Value = '3.615.35'
$ws.Cells.Item(24, 5).Value = '  -1.22%  '
$ws.Cells.Item(25, 4).Value = '''73.02'
$ws.Cells.Item(25, 5).Value = '  -1.26%  '
$ws.Cells.Item(26, 4).Value = '''1.00'
$ws.Cells.Item(26, 5).Value = '  -0.09%  '
$ws.Cells.Item(27, 5).Value = '  +1.44%  '
$ws.Cells.Item(28, 5).Value = '  -1.23%  '
$ws.Cells.Item(29, 4).Value = '''0.999'
$ws.Cells.Item(29, 5).Value = '  -0.44%  '
$ws.Cells.Item(30, 4).Value = '''7.40'
$ws.Cells.Item(30, 5).Value = '  -3.27%  '
$ws.Cells.Item(31, 4).Value = '''8.16'
$ws.Cells.Item(31, 5).Value = '  -1.53%  '
$ws.Cells.Item(32, 5).Value = '  -3.25%  '
$ws.Cells.Item(33, 4).Value = '3.484.72'
$ws.Cells.Item(33, 5).Value = '  -1.02%  '
$ws.Cells.Item(34, 5).Value = '  -0.02%  '
$ws.Cells.Item(35, 4).Value = '''23.26'
$ws.Cells.Item(35, 5).Value = '  -3.73%  '
$ws.Cells.Item(36, 5).Value = '  -0.59%  '
$ws.Cells.Item(37, 5).Value = '  -0.06%  '
$ws.Cells.Item(38, 5).Value = '  -0.12%  '
$ws.Cells.Item(39, 5).Value = '  -1.84%  '
$ws.Cells.Item(40, 4).Value = '''159.73'
$ws.Cells.Item(40, 5).Value = '  -5.18%  '
$ws.Cells.Item(41, 5).Value = '  -3.92%  '
$ws.Cells.Item(42, 2).Value = 'Mantle'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Cells.Item(42, 4).Value = '''0.808'
$ws.Cells.Item(42, 5).Value = '  -1.66%  '
$ws.Cells.Item(43, 2).Value = 'EnergySwap'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(43, 4).Value = '''25.89'
$ws.Cells.Item(43, 5).Value = '  +0.49%  '
$ws.Cells.Item(44, 5).Value = '  +0.05%  '
$ws.Cells.Item(45, 4).Value = '''41.58'
$ws.Cells.Item(45, 5).Value = '  -2.51%  '
$ws.Cells.Item(46, 4).Value = '''1.19'
$ws.Cells.Item(46, 5).Value = '  -4.00%  '
$ws.Cells.Item(47, 5).Value = '  -2.19%  '
$ws.Cells.Item(48, 5).Value = '  -2.65%  '
$ws.Cells.Item(49, 4).Value = '2.413.15'
$ws.Cells.Item(49, 5).Value = '  +1.38%  '
$ws.Cells.Item(50, 5).Value = '  -1.68%  '
$ws.Cells.Item(51, 4).Value = '''0.882'
$ws.Cells.Item(51, 5).Value = '  -0.67%  '
